$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (rows 5 and 6) so only 3 data rows remain
$ws.Rows("5:6").Delete()

# Build a helper cell holding the literal text "2023-07-18" (via a TEXT()
# formula) so that pasting its value into the date column stores a plain
# shared string instead of having Excel auto-convert the date-like text
# into a date serial number.
$ws.Range("H1").Formula = '=TEXT(DATE(2023,7,18),"yyyy-mm-dd")'
$ws.Range("H1").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("A3").PasteSpecial(-4163)
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("H1").Clear()

# Row 2: Buy ODAS.IS
$ws.Range("B2").Value = "Buy"
$ws.Range("C2").Value = "ODAS.IS"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 34
$ws.Range("F2").Value = 102

# Row 3: Buy EREGL.IS
$ws.Range("B3").Value = "Buy"
$ws.Range("C3").Value = "EREGL.IS"
$ws.Range("D3").Value = 40
$ws.Range("E3").Value = 34
$ws.Range("F3").Value = 1360

# Row 4: Sell SISE.IS
$ws.Range("B4").Value = "Sell"
$ws.Range("C4").Value = "SISE.IS"
$ws.Range("D4").Value = 34
$ws.Range("E4").Value = -3
$ws.Range("F4").Value = -102
